# Apply updated cryptocurrency price/volume figures (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper (scratch) cell used to build a pure-text value via a formula, then
# paste-special as a value so the target cell keeps its original (default) style
# and Excel does not "helpfully" reinterpret numeric-looking text (e.g. "223.51")
# as a floating point number.
$helper = $ws.Range("ZZ1")

function Set-TextValue {
    param($targetAddress, $text)
    $helper.Formula = '="' + $text + '"'
    $helper.Copy() | Out-Null
    $ws.Range($targetAddress).PasteSpecial(-4163) | Out-Null
    $helper.ClearContents() | Out-Null
}

Set-TextValue "D2" '27.317.70'
Set-TextValue "E2" '  -2.58%  '
Set-TextValue "D3" '1.708.83'
Set-TextValue "E3" '  -1.94%  '
Set-TextValue "E4" '  +0.18%  '
Set-TextValue "D5" '223.51'
Set-TextValue "E5" '  -2.81%  '
Set-TextValue "D6" '0.5310'
Set-TextValue "E6" '  -2.83%  '
Set-TextValue "E7" '  +0.17%  '
Set-TextValue "E8" '  -4.83%  '
Set-TextValue "D9" '0.06600'
Set-TextValue "E9" '  -2.06%  '
Set-TextValue "D10" '20.89'
Set-TextValue "E10" '  -4.79%  '
Set-TextValue "D11" '0.07654'
Set-TextValue "E11" '  -1.76%  '
Set-TextValue "E12" '  -2.96%  '
Set-TextValue "D13" '1.723.79'
Set-TextValue "E13" '  -0.82%  '
Set-TextValue "D14" '1.942.71'
Set-TextValue "E14" '  -1.94%  '
Set-TextValue "D15" '0.5742'
Set-TextValue "E15" '  -4.79%  '
Set-TextValue "E16" '  -3.09%  '
Set-TextValue "D17" '67.61'
Set-TextValue "E17" '  -3.34%  '
Set-TextValue "D18" '27.297.82'
Set-TextValue "E18" '  -2.52%  '
Set-TextValue "D19" '216.61'
Set-TextValue "E19" '  -5.51%  '
Set-TextValue "E20" '  +0.11%  '
Set-TextValue "D21" '4.678'
Set-TextValue "E21" '  -3.54%  '
Set-TextValue "E22" '  -5.04%  '
Set-TextValue "D23" '5.982'
Set-TextValue "E23" '  -4.87%  '
Set-TextValue "E24" '  +0.11%  '
Set-TextValue "D25" '142.45'
Set-TextValue "E25" '  -3.22%  '
Set-TextValue "D26" '1.747'
Set-TextValue "E26" '  +7.99%  '
Set-TextValue "D27" '0.1217'
Set-TextValue "E27" '  -2.86%  '
Set-TextValue "D28" '7.264'
Set-TextValue "E28" '  -2.92%  '
Set-TextValue "D29" '16.34'
Set-TextValue "E29" '  -5.08%  '
Set-TextValue "D30" '0.05377'
Set-TextValue "E30" '  -5.00%  '
Set-TextValue "D31" '1.292'
Set-TextValue "E31" '  -2.15%  '
Set-TextValue "E32" '  -5.69%  '
Set-TextValue "D33" '3.424'
Set-TextValue "E33" '  -3.79%  '
Set-TextValue "D34" '1.640'
Set-TextValue "E34" '  -1.37%  '
Set-TextValue "D35" '2.879'
Set-TextValue "E35" '  +0.59%  '
Set-TextValue "D36" '2.422'
Set-TextValue "E36" '  -1.21%  '
Set-TextValue "D37" '0.9488'
Set-TextValue "E37" '  -3.92%  '
Set-TextValue "D38" '0.5860'
Set-TextValue "E38" '  -1.69%  '
Set-TextValue "D39" '0.01634'
Set-TextValue "E39" '  -2.89%  '
Set-TextValue "D40" '5.865'
Set-TextValue "E40" '  -2.27%  '
Set-TextValue "E41" '  +0.11%  '
Set-TextValue "D42" '1.041.75'
Set-TextValue "E42" '  -0.80%  '
Set-TextValue "D43" '0.8406'
Set-TextValue "E43" '  -0.91%  '
Set-TextValue "D44" '101.06'
Set-TextValue "D45" '1.849.55'
Set-TextValue "E45" '  -1.85%  '
Set-TextValue "E46" '  -1.81%  '
Set-TextValue "D47" '58.14'
Set-TextValue "E47" '  -3.78%  '
Set-TextValue "D48" '0.4502'
Set-TextValue "E48" '  +1.83%  '
Set-TextValue "D49" '1.008'
Set-TextValue "E49" '  -0.41%  '
Set-TextValue "D50" '8.072'
Set-TextValue "E50" '  -3.18%  '
Set-TextValue "D51" '0.06522'
Set-TextValue "E51" '  +10.01%  '

$excel.CutCopyMode = 0
